$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old "a" / "R6..." entries in row 16 with the new R6 redirect text,
# and fill in the N / Minor columns that were missing before.
$ws.Range("A16").Value = "R6 The application shall redirect users to find people page"
$ws.Range("B16").Value = "N"
$ws.Range("C16").Value = "Minor"

# Add new rows 18, 20, 22, 24, 26 (odd rows 17/19/21/23/25 stay blank, matching
# the existing pattern of skipped row 13).
$ws.Range("A18").Value = "R7 The application shall redirect users to news flash page"
$ws.Range("B18").Value = "N"
$ws.Range("C18").Value = "Minor"

$ws.Range("A20").Value = "R8 The application shall redirect users to mymav webpage"
$ws.Range("B20").Value = "N"
$ws.Range("C20").Value = "Minor"

$ws.Range("A22").Value = "R9 The application shall redirect users to uta homepage"
$ws.Range("B22").Value = "R"
$ws.Range("C22").Value = "Minor"

$ws.Range("A24").Value = "R10 The application shall redirect users to alumini page"
$ws.Range("B24").Value = "D"
$ws.Range("C24").Value = "Minor"

$ws.Range("A26").Value = "R 11 The appliication shall provide link to uta email"
$ws.Range("B26").Value = "R"
$ws.Range("C26").Value = "Minor"

$ws.Range("C26").Select()
